$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 46: KKR vs KXI (row 55) - enter the six players' raw scores.
# Formulas in D/G/J/M/P/S recompute automatically via RANK/VLOOKUP.
$ws.Range("E55").Value = 40
$ws.Range("H55").Value = 100
$ws.Range("K55").Value = 20
$ws.Range("N55").Value = 80
$ws.Range("Q55").Value = 60
$ws.Range("T55").Value = 0
